# Weekly update for "Fruta, Terminal Hortofrutícola Agro Chillán - Plátano":
# a new week of price observations (2 rows: "Pintón" and "Primera Pintón")
# is inserted at the top of the data block (row 1073), pushing every
# existing data row down by two. The data range therefore grows from
# A1:T1142 to A1:T1144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block; Excel shifts rows
# 1073:1142 down to 1075:1144 and keeps their contents/formatting intact.
$ws.Rows("1073:1074").Insert()

# New row 1073 - "Pintón" quality, week of 2023-12-05 (serial 45265).
$ws.Range("A1073").Value = 7
$ws.Range("B1073").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C1073").Value = 'Ñuble'
$ws.Range("D1073").Value = 45265
$ws.Range("E1073").Value = 16
$ws.Range("F1073").Value = 'Fruta'
$ws.Range("G1073").Value = 100108
$ws.Range("H1073").Value = 'Tropicales y subtropicales'
$ws.Range("I1073").Value = 100108006
$ws.Range("J1073").Value = 'Plátano'
$ws.Range("K1073").Value = 'Sin especificar'
$ws.Range("L1073").Value = 'Pintón'
$ws.Range("M1073").Value = 250
$ws.Range("N1073").Value = 25000
$ws.Range("O1073").Value = 25000
$ws.Range("P1073").Value = 25000
$ws.Range("Q1073").Value = '$/caja 20 kilos'
$ws.Range("R1073").Value = 'Ecuador'
$ws.Range("S1073").Value = 1250
$ws.Range("T1073").Value = 20

# New row 1074 - "Primera Pintón" quality, same week.
$ws.Range("A1074").Value = 7
$ws.Range("B1074").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C1074").Value = 'Ñuble'
$ws.Range("D1074").Value = 45265
$ws.Range("E1074").Value = 16
$ws.Range("F1074").Value = 'Fruta'
$ws.Range("G1074").Value = 100108
$ws.Range("H1074").Value = 'Tropicales y subtropicales'
$ws.Range("I1074").Value = 100108006
$ws.Range("J1074").Value = 'Plátano'
$ws.Range("K1074").Value = 'Sin especificar'
$ws.Range("L1074").Value = 'Primera Pintón'
$ws.Range("M1074").Value = 200
$ws.Range("N1074").Value = 26000
$ws.Range("O1074").Value = 26000
$ws.Range("P1074").Value = 26000
$ws.Range("Q1074").Value = '$/caja 20 kilos'
$ws.Range("R1074").Value = 'Ecuador'
$ws.Range("S1074").Value = 1300
$ws.Range("T1074").Value = 20
